$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

function Copy-RowTo($srcRow, $dstRow) {
    $src = $ws.Range("A" + $srcRow + ":AY" + $srcRow)
    $dst = $ws.Range("A" + $dstRow + ":AY" + $dstRow)
    # Clear the destination first: PasteSpecial leaves a destination cell
    # untouched whenever the matching source cell is blank/absent, so
    # without this, stale values (e.g. a comment that shouldn't carry
    # over) would survive the "move".
    $dst.ClearContents() | Out-Null
    $src.Copy() | Out-Null
    $dst.PasteSpecial($xlPasteValues) | Out-Null
}

# The data rows 2-6 (columns A:AY) get cyclically re-shuffled so that each
# row ends up holding what used to be a different row's data:
#   new row 2 <= old row 6
#   new row 3 <= old row 5
#   new row 4 <= old row 2
#   new row 5 <= old row 3
#   new row 6 <= old row 4
#
# Equivalently, thinking of where each OLD row's data travels TO:
#   old row 2 -> new row 4   \
#   old row 4 -> new row 6    |  3-cycle (2 -> 4 -> 6 -> 2)
#   old row 6 -> new row 2   /
#   old row 3 -> new row 5   \  2-cycle / swap (3 <-> 5)
#   old row 5 -> new row 3   /
#
# A scratch row well outside the used range holds one row's original
# values while the cycle is resolved, using Copy/PasteSpecial(values) so
# text that looks like a date (e.g. "2023-09-01") is carried over as
# literal text instead of being reinterpreted as a date serial number.
$scratchRow = 200

# --- 3-cycle: 2 -> 4 -> 6 -> 2 ---
Copy-RowTo 2 $scratchRow   # stash old row 2
Copy-RowTo 6 2             # old row 6 -> row 2
Copy-RowTo 4 6             # old row 4 -> row 6
Copy-RowTo $scratchRow 4   # stashed old row 2 -> row 4

# --- 2-cycle: 3 <-> 5 ---
Copy-RowTo 3 $scratchRow   # stash old row 3
Copy-RowTo 5 3             # old row 5 -> row 3
Copy-RowTo $scratchRow 5   # stashed old row 3 -> row 5

# Clear the scratch row so it doesn't linger in the saved workbook.
$ws.Range("A" + $scratchRow + ":AY" + $scratchRow).ClearContents() | Out-Null

# Row 7: only the Taxonsorteringsordning (column B) value changes.
$ws.Range("B7").Value2 = 89557
